# [feature] Read Excel && Write File Init
#
# Rebuilds the Sheet1 header/type rows with the new column layout
# (ID/name/age/hight/sex/location over int/string/int/float/bool/vector3),
# drops the old rotation/weight/TestValue columns (G:I), and closes the
# gap in the ID column so rows 4-8 hold 100001..100005 contiguously.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old G:I columns (rotation / Height / TestValue row) entirely -
# the new layout only uses columns A:F.
[void]$ws.Range("G1:I20").ClearContents()

# Row 2 = field names, row 3 = field types (column by column so duplicate
# values like "int" line up the way the source workbook built them).
$ws.Range("A2").Value = "ID"
$ws.Range("A3").Value = "int"

$ws.Range("B2").Value = "name"
$ws.Range("B3").Value = "string"

$ws.Range("C2").Value = "age"
$ws.Range("C3").Value = "int"

$ws.Range("D2").Value = "hight"
$ws.Range("D3").Value = "float"

$ws.Range("E2").Value = "sex"
$ws.Range("E3").Value = "bool"

$ws.Range("F2").Value = "location"
$ws.Range("F3").Value = "vector3"

# Data rows: close the gap that used to leave row 6 blank - IDs now run
# contiguously from row 4 through row 8.
$ws.Range("A4").Value = 100001
$ws.Range("A5").Value = 100002
$ws.Range("A6").Value = 100003
$ws.Range("A7").Value = 100004
$ws.Range("A8").Value = 100005

# The old row 9 (previously holding 100005) is no longer part of the used
# range - clear any leftovers from it.
[void]$ws.Range("A9:I9").ClearContents()

# Match the workbook's saved selection state.
[void]$ws.Range("F7").Select()
